# Vocabulary workbook update: new .ttl generated from Google sheet.
# Updates skos:broader (column F) values for several existing rows and
# appends 20 new concept rows (214-233) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update skos:broader (column F) on existing rows ------------------
$ws.Range("F88").Value  = "covid19:10108"
$ws.Range("F89").Value  = "covid19:10108"
$ws.Range("F90").Value  = "covid19:10108"
$ws.Range("F91").Value  = "covid19:10132"
$ws.Range("F92").Value  = "covid19:10108"
$ws.Range("F93").Value  = "covid19:10109"
$ws.Range("F97").Value  = "covid19:10108"
$ws.Range("F99").Value  = "covid19:10108"
$ws.Range("F100").Value = "covid19:10108"
$ws.Range("F189").Value = "id-amr:10162, gen:10005"

# --- Append new rows 214-233 -------------------------------------------
$newRows = @(
    @(214, "covid19:10109", "health record data", "covid19:10108"),
    @(215, "id-amr:10190", "physiological - biochemical measurements", "covid19:10112"),
    @(216, "id-amr:10191", "sequencing - genotyping data", "covid19:10244"),
    @(217, "id-amr:10192", "clinical records", "covid19:10108"),
    @(218, "id-amr:10193", "follow-up records", "covid19:10108"),
    @(219, "id-amr:10194", "monitoring data ", "covid19:10108"),
    @(220, "id-amr:10195", "surveillance data", "covid19:10108"),
    @(221, "id-amr:10196", "patient ethnicity data", "covid19:10262"),
    @(222, "id-amr:10197", "quality data", "gen:10005"),
    @(223, "id-amr:10198", "experimental data ", "gen:10005"),
    @(224, "id-amr:10199", "biospecimen processing", "gen:10013"),
    @(225, "id-amr:10200", "data analysis", "gen:10013"),
    @(226, "id-amr:10201", "data integration", "gen:10013"),
    @(227, "id-amr:10202", "digital imaging", "gen:10013"),
    @(228, "id-amr:10203", "FAIRification service", "gen:10013"),
    @(229, "id-amr:10204", "providing reference materials", "gen:10013"),
    @(230, "id-amr:10205", "quality control", "gen:10013"),
    @(231, "id-amr:10206", "regulatory", "gen:10013"),
    @(232, "id-amr:10207", "training", "gen:10013"),
    @(233, "id-amr:10208", "communication", "gen:10013")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("F$rowNum").Value = $r[3]
}
